# Auto-generated edit script: updates crypto price/volume table (Sheet1)
# to match the scraped values from the "Thu May  9 03:40:59 UTC 2024" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.605.34'
$ws.Range('E2').Value = '  -1.78%  '

# Row 3
$ws.Range('D3').Value = '2.997.20'
$ws.Range('E3').Value = '  -1.00%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.06%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.28%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.522'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.65%  '

# Row 9
$ws.Range('D9').Value = '2.997.90'
$ws.Range('E9').Value = '  -0.99%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.48%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.93'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.04%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.463'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.23%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000229'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.99%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.81%  '

# Row 15
$ws.Range('E15').Value = '  +2.06%  '

# Row 16
$ws.Range('D16').Value = '3.486.27'
$ws.Range('E16').Value = '  -1.17%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.99%  '

# Row 18
$ws.Range('D18').Value = '61.508.27'
$ws.Range('E18').Value = '  -1.87%  '

# Row 19
$ws.Range('D19').Value = '3.007.49'
$ws.Range('E19').Value = '  -0.61%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '454.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.01%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.688'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.63%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.54%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.70%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.09%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.53%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.10%  '

# Row 28
$ws.Range('E28').Value = '  +0.19%  '

# Row 29
$ws.Range('E29').Value = '  +2.01%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.11%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.64%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.34%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.68'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.84%  '

# Row 34
$ws.Range('E34').Value = '  +0.30%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0820'
$ws.Range('E35').Value = '  +2.89%  '

# Row 36
$ws.Range('E36').Value = '  -1.68%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.50%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.40%  '

# Row 39
$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.81%  '

# Row 40
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.25%  '

# Row 41
$ws.Range('E41').Value = '  +7.74%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.87'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.55%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '398.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.15%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0353'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.83%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '38.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.10%  '

# Row 46
$ws.Range('E46').Value = '  -5.14%  '

# Row 47
$ws.Range('D47').Value = '2.720.39'
$ws.Range('E47').Value = '  -2.94%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.58%  '

# Row 50
$ws.Range('E50').Value = '  -0.72%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.13%  '
